# The "9.7国乙同人ONLY（取消）" event (row 2) was cancelled and has been
# dropped from the listing. Remove it from every sheet that lists it
# ("展览" and "全部类型" — the other two sheets, "演出" and "本地生活", never
# contained this event). All following rows shift up by one, and a handful
# of "想去人数" (want-to-go) counts were refreshed to their latest totals.

$wb = $excel.ActiveWorkbook

function Update-SinaExpoSheet($ws) {
    # Remove the cancelled event row; everything below shifts up one row.
    $ws.Rows.Item(2).Delete()

    # Column A holds literal sequence numbers (1, 2, 3, …) that Excel does
    # not renumber automatically on a row delete, so restamp them to match
    # the new row order.
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value2 = $r - 1
    }

    # Refresh the "想去人数" (want-to-go count) figures that changed since
    # the last scrape, matched by the event's unique bilibili show id in
    # column H.
    $updates = @{
        "id=91385" = 206
        "id=90762" = 3542
        "id=91043" = 361
        "id=87820" = 428
    }

    for ($r = 2; $r -le $lastRow; $r++) {
        $link = $ws.Cells.Item($r, 8).Value2
        foreach ($key in $updates.Keys) {
            if ($link -like "*$key*") {
                $ws.Cells.Item($r, 6).Value2 = $updates[$key]
            }
        }
    }
}

foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "展览" -or $ws.Name -eq "全部类型") {
        Update-SinaExpoSheet $ws
    }
}
